$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '41.656.81'
$cell.Style = "Normal"
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = '  +0.31%  '
$cell.Style = "Normal"
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '2.467.47'
$cell.Style = "Normal"
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = '  -0.69%  '
$cell.Style = "Normal"
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = '  +0.23%  '
$cell.Style = "Normal"
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '316.60'
$cell.Style = "Normal"
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = '  +0.43%  '
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '92.71'
$cell.Style = "Normal"
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = '  -0.52%  '
$cell.Style = "Normal"
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.552'
$cell.Style = "Normal"
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = '  +1.84%  '
$cell.Style = "Normal"
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = '  +0.04%  '
$cell.Style = "Normal"
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.517'
$cell.Style = "Normal"
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = '  +2.76%  '
$cell.Style = "Normal"
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.0900'
$cell.Style = "Normal"
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = '  +14.79%  '
$cell.Style = "Normal"
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '32.69'
$cell.Style = "Normal"
$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = '  -0.17%  '
$cell.Style = "Normal"
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = '  -0.45%  '
$cell.Style = "Normal"
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '2.848.74'
$cell.Style = "Normal"
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = '  -0.61%  '
$cell.Style = "Normal"
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = '  +0.98%  '
$cell.Style = "Normal"
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '15.63'
$cell.Style = "Normal"
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = '  -1.71%  '
$cell.Style = "Normal"
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '2.474.98'
$cell.Style = "Normal"
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = '  -2.50%  '
$cell.Style = "Normal"
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = '  +3.63%  '
$cell.Style = "Normal"
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '41.627.85'
$cell.Style = "Normal"
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = '  +0.26%  '
$cell.Style = "Normal"
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0971'
$cell.Style = "Normal"
$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = '  +4.37%  '
$cell.Style = "Normal"
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '6.45'
$cell.Style = "Normal"
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = '  +1.68%  '
$cell.Style = "Normal"
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '71.37'
$cell.Style = "Normal"
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = '  -0.03%  '
$cell.Style = "Normal"
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '11.44'
$cell.Style = "Normal"
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = '  +2.01%  '
$cell.Style = "Normal"
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '240.64'
$cell.Style = "Normal"
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = '  +1.45%  '
$cell.Style = "Normal"
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = '  +0.37%  '
$cell.Style = "Normal"
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = '  +0.42%  '
$cell.Style = "Normal"
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = '  -0.06%  '
$cell.Style = "Normal"
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '24.77'
$cell.Style = "Normal"
$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = '  -1.29%  '
$cell.Style = "Normal"
$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = '  +3.19%  '
$cell.Style = "Normal"
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '9.81'
$cell.Style = "Normal"
$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = '  +1.22%  '
$cell.Style = "Normal"
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '35.27'
$cell.Style = "Normal"
$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = '  -2.21%  '
$cell.Style = "Normal"
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '155.99'
$cell.Style = "Normal"
$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = '  -0.90%  '
$cell.Style = "Normal"
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '5.49'
$cell.Style = "Normal"
$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = '  +0.62%  '
$cell.Style = "Normal"
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '2.57'
$cell.Style = "Normal"
$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = '  -1.01%  '
$cell.Style = "Normal"
$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = '  +1.49%  '
$cell.Style = "Normal"
$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = '  +1.41%  '
$cell.Style = "Normal"
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '17.59'
$cell.Style = "Normal"
$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = '  -0.76%  '
$cell.Style = "Normal"
$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = '  -1.25%  '
$cell.Style = "Normal"
$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = '  +1.48%  '
$cell.Style = "Normal"
$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = '  -1.95%  '
$cell.Style = "Normal"
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = '  -3.50%  '
$cell.Style = "Normal"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '3.96'
$cell.Style = "Normal"
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = '  -3.48%  '
$cell.Style = "Normal"
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = '  -0.06%  '
$cell.Style = "Normal"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '1.978.85'
$cell.Style = "Normal"
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = '  +0.48%  '
$cell.Style = "Normal"
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = '  +0.35%  '
$cell.Style = "Normal"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '18.84'
$cell.Style = "Normal"
$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = '  -4.09%  '
$cell.Style = "Normal"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '2.98'
$cell.Style = "Normal"
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = '  +0.35%  '
$cell.Style = "Normal"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '9.11'
$cell.Style = "Normal"
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = '  +1.99%  '
$cell.Style = "Normal"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '2.703.11'
$cell.Style = "Normal"
$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = '  -0.77%  '
$cell.Style = "Normal"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '96.97'
$cell.Style = "Normal"
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = '  -0.34%  '
$cell.Style = "Normal"
$cell = $ws.Range("B50")
$cell.NumberFormat = "@"
$cell.Value = 'BitcoinSV'
$cell.Style = "Normal"
$cell = $ws.Range("C50")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$cell.Style = "Normal"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '74.21'
$cell.Style = "Normal"
$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = '  +1.74%  '
$cell.Style = "Normal"
$cell = $ws.Range("B51")
$cell.NumberFormat = "@"
$cell.Value = 'ordi'
$cell.Style = "Normal"
$cell = $ws.Range("C51")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$cell.Style = "Normal"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '66.74'
$cell.Style = "Normal"
$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = '  -1.40%  '
$cell.Style = "Normal"
